# Generate Report for Handback
# Update the timestamp strings recorded in the handback status workbook.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G4").Value = "2016-08-22 14:51:47"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H4").Value = "2016-08-22 14:51:42"
$wsZhCn.Range("K4").Value = "2016-08-22 14:52:04"

$wsDeDe = $wb.Worksheets.Item("de-de")
# de-de!H4 shares the same "Latest HO Xliff Generate Date" value as Overview!G4
$wsDeDe.Range("H4").Value = "2016-08-22 14:51:47"
$wsDeDe.Range("K4").Value = "2016-08-22 14:52:25"
